# Auto-generated update of DAMSLTag (col I) and DialogAct (col J) values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=25; I="sv"; J="Statement-opinion"}
    @{Row=29; I="ba"; J="Appreciation"}
    @{Row=37; I="ba"; J="Appreciation"}
    @{Row=43; I="sv"; J="Statement-opinion"}
    @{Row=47; I="sv"; J="Statement-opinion"}
    @{Row=63; I="sv"; J="Statement-opinion"}
    @{Row=70; I="sv"; J="Statement-opinion"}
    @{Row=72; I="sd"; J="Statement-non-opinion"}
    @{Row=74; I="sd"; J="Statement-non-opinion"}
    @{Row=78; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=80; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=89; I="sd"; J="Statement-non-opinion"}
    @{Row=98; I="sd"; J="Statement-non-opinion"}
    @{Row=107; I="sv"; J="Statement-opinion"}
    @{Row=131; I="sd"; J="Statement-non-opinion"}
    @{Row=142; I="sv"; J="Statement-opinion"}
    @{Row=144; I="sd"; J="Statement-non-opinion"}
    @{Row=146; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=150; I="sd"; J="Statement-non-opinion"}
    @{Row=156; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=167; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=170; I="sd"; J="Statement-non-opinion"}
    @{Row=183; I="sv"; J="Statement-opinion"}
    @{Row=185; I="ba"; J="Appreciation"}
    @{Row=193; I="aa"; J="Agree/Accept"}
    @{Row=197; I="ba"; J="Appreciation"}
    @{Row=201; I="sd"; J="Statement-non-opinion"}
    @{Row=202; I="sv"; J="Statement-opinion"}
    @{Row=208; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=221; I="sv"; J="Statement-opinion"}
    @{Row=224; I="sv"; J="Statement-opinion"}
    @{Row=233; I="%"; J="Uninterpretable"}
    @{Row=234; I="qy"; J="Yes-No-Question"}
    @{Row=236; I="aa"; J="Agree/Accept"}
    @{Row=244; I="aa"; J="Agree/Accept"}
    @{Row=245; I="aa"; J="Agree/Accept"}
    @{Row=248; I="sd"; J="Statement-non-opinion"}
    @{Row=261; I="sd"; J="Statement-non-opinion"}
    @{Row=274; I="sd"; J="Statement-non-opinion"}
    @{Row=275; I="sv"; J="Statement-opinion"}
    @{Row=278; I="sv"; J="Statement-opinion"}
    @{Row=279; I="sd"; J="Statement-non-opinion"}
    @{Row=294; I="aa"; J="Agree/Accept"}
    @{Row=305; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=307; I="sv"; J="Statement-opinion"}
    @{Row=327; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=352; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=358; I="sd"; J="Statement-non-opinion"}
    @{Row=382; I="aa"; J="Agree/Accept"}
    @{Row=385; I="aa"; J="Agree/Accept"}
    @{Row=394; I="sv"; J="Statement-opinion"}
    @{Row=402; I="sd"; J="Statement-non-opinion"}
    @{Row=410; I="sv"; J="Statement-opinion"}
    @{Row=419; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=428; I="sd"; J="Statement-non-opinion"}
    @{Row=432; I="sv"; J="Statement-opinion"}
    @{Row=437; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=440; I="sv"; J="Statement-opinion"}
    @{Row=462; I="sd"; J="Statement-non-opinion"}
    @{Row=483; I="sv"; J="Statement-opinion"}
    @{Row=493; I="sv"; J="Statement-opinion"}
    @{Row=495; I="sv"; J="Statement-opinion"}
    @{Row=499; I="ba"; J="Appreciation"}
    @{Row=503; I="ba"; J="Appreciation"}
    @{Row=504; I="ba"; J="Appreciation"}
    @{Row=506; I="sd"; J="Statement-non-opinion"}
    @{Row=520; I="sd"; J="Statement-non-opinion"}
    @{Row=526; I="sv"; J="Statement-opinion"}
    @{Row=542; I="sd"; J="Statement-non-opinion"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
